$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values per the transaction edit
$ws.Range("B2").Value = 45324
$ws.Range("D2").Value = "groceries"
$ws.Range("E2").Value = 150
$ws.Range("F2").Value = "February groceries note"
